$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows.Item(24).Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
